# 秦皇岛动力煤价格_月度数据 — append the new 2025/12/31 monthly data point
# (the sheet is sorted most-recent-first, so this is a row insert at row 2)
# and refresh the forecast/actual values that moved with the new data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the first data row, shifting history down.
$ws.Rows("2:2").Insert()

# The inserted row inherits the header row's bold/border/centered style;
# reset it back to the plain style the rest of the data rows use.
$ws.Range("A2:C2").Style = "Normal"

# New data point for 2025/12/31 (真实值/actual not published yet, so B2
# is left blank just like the inserted row's default empty cell).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025/12/31"
$ws.Range("C2").Value = 770.8

# Refresh 真实值 (actual) for 2025/10/31, now that the month has closed.
$ws.Range("B4").Value = 740.9

# Refresh 预测值 (forecast) values across history with the latest pull.
$ws.Range("C3").Value = 771.4
$ws.Range("C4").Value = 785.3
$ws.Range("C5").Value = 785.5
$ws.Range("C6").Value = 762.7
$ws.Range("C7").Value = 712.2
$ws.Range("C8").Value = 683.6
$ws.Range("C11").Value = 762.5
$ws.Range("C12").Value = 775.5
$ws.Range("C13").Value = 766.2
$ws.Range("C14").Value = 774.8
$ws.Range("C15").Value = 805.7
$ws.Range("C16").Value = 839.9
$ws.Range("C17").Value = 859.2
$ws.Range("C18").Value = 846.2
$ws.Range("C19").Value = 820.3
$ws.Range("C20").Value = 821.4
$ws.Range("C21").Value = 822
$ws.Range("C22").Value = 838.8
$ws.Range("C23").Value = 860.6
$ws.Range("C24").Value = 868.1
$ws.Range("C25").Value = 866.6
